# Update odds values on Sheet1 to reflect the refreshed FlashScore snapshot.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4 (Patriotas - Santa Fe)
$ws.Range("G4").Value = 4.33
$ws.Range("H4").Value = 3.1
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 5
$ws.Range("K4").Value = 1.91
$ws.Range("L4").Value = 2.75
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.6
$ws.Range("R4").Value = 1.48
$ws.Range("S4").Value = 1.57
$ws.Range("T4").Value = 2.25
$ws.Range("U4").Value = 2.25
$ws.Range("V4").Value = 1.57
$ws.Range("W4").Value = 9
$ws.Range("X4").Value = 19
$ws.Range("Z4").Value = 51
$ws.Range("AD4").Value = 6
$ws.Range("AH4").Value = 8
$ws.Range("AI4").Value = 9.5
$ws.Range("AJ4").Value = 17
$ws.Range("AK4").Value = 21
$ws.Range("AN4").Value = 6
$ws.Range("AO4").Value = 26
$ws.Range("AQ4").Value = 101
$ws.Range("AT4").Value = 2.25
$ws.Range("AW4").Value = 3.75
$ws.Range("AX4").Value = 12
$ws.Range("AZ4").Value = 41
$ws.Range("BB4").Value = 251

# Row 5 (Once Caldas - Dep. Pasto)
$ws.Range("O5").Value = 1.5
$ws.Range("P5").Value = 2.5

# Row 6 (America De Cali - Ind. Medellin)
$ws.Range("O6").Value = 1.44
$ws.Range("P6").Value = 2.63
$ws.Range("Q6").Value = 2.4
$ws.Range("R6").Value = 1.53

# Row 9 (Bohemians - St. Patricks)
$ws.Range("G9").Value = 2.75
$ws.Range("I9").Value = 2.25
$ws.Range("J9").Value = 3.4
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 1.04
$ws.Range("N9").Value = 12
$ws.Range("Q9").Value = 1.85
$ws.Range("R9").Value = 2
$ws.Range("X9").Value = 15
$ws.Range("Y9").Value = 11
$ws.Range("Z9").Value = 29
$ws.Range("AB9").Value = 29
$ws.Range("AG9").Value = 9
$ws.Range("AH9").Value = 12
$ws.Range("AI9").Value = 9.5
$ws.Range("AJ9").Value = 23
$ws.Range("AN9").Value = 5
$ws.Range("AP9").Value = 23
$ws.Range("AQ9").Value = 51
$ws.Range("AW9").Value = 4.5
$ws.Range("BB9").Value = 126

# Row 10 (Gijon - Castellon)
$ws.Range("H10").Value = 3.5
$ws.Range("I10").Value = 3.1
$ws.Range("K10").Value = 2.25
$ws.Range("L10").Value = 3.5
$ws.Range("O10").Value = 1.25
$ws.Range("P10").Value = 3.75
$ws.Range("Q10").Value = 1.83
$ws.Range("R10").Value = 2.03
$ws.Range("U10").Value = 1.62
$ws.Range("V10").Value = 2.2
$ws.Range("W10").Value = 9
$ws.Range("AE10").Value = 13
$ws.Range("AG10").Value = 11
$ws.Range("AK10").Value = 23
$ws.Range("AL10").Value = 29
$ws.Range("AM10").Value = 151
$ws.Range("AN10").Value = 4.33
$ws.Range("AP10").Value = 21
$ws.Range("AS10").Value = 126
$ws.Range("AV10").Value = 51
$ws.Range("AX10").Value = 17
$ws.Range("AY10").Value = 23
$ws.Range("BA10").Value = 67
$ws.Range("BB10").Value = 151

$wb.Save()
